$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 1: extend the sequence with two more columns (P1=14, Q1=15),
# using the same formatting as the existing header cell O1.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Data rows 2-25: columns I and K swap their values, columns M and O swap
# their values, and two new columns P and Q are appended (always 2).
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2    # I
    $ws.Cells.Item($r, 11).Value = 1   # K
    $ws.Cells.Item($r, 13).Value = 2   # M
    $ws.Cells.Item($r, 15).Value = 1   # O
    $ws.Cells.Item($r, 16).Value = 2   # P
    $ws.Cells.Item($r, 17).Value = 2   # Q
}
